$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.597.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.684.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.39%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5334"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +5.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06420"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.686.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.505"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5606"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8413"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.646.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.43%  "
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.790"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "196.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.35%  "
$ws.Range("E21").Value = "  +4.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.360"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.99%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1279"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.472"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.436"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06172"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.278"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.611"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.467"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.698"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.010"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.45%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.424"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.796"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5746"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01645"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.038"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.071.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8624"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.834.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₈108"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.163"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.26%  "
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.071"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.52%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4243"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.63%  "
